# Applies the "Removed double type and re-task float" edit:
#  - "double" type identifier row (A10) is removed/repurposed as "reserved"
#  - shared strings "an integer?" and "double" are dropped entirely
#  - "[double]" becomes "[float]" (K18)
#  - K15 ("integer") -> "long", K16 ("an integer?") -> "integer" (types retask)
#  - D18 flips from 1 to 0 (geo row no longer includes that flag bit),
#    which recalculates J18's formula result from 142 to 138
#  - selection moves from C5 to C18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "reserved"
$ws.Range("K15").Value = "long"
$ws.Range("K16").Value = "integer"
$ws.Range("K18").Value = "[float]"

$ws.Range("D18").Value = 0

$ws.Range("C18").Select()

$wb.Save()
